# Apply text replacements to the two-digit division answer table.
# Each (old -> new) pair corresponds to one <w:t> run in document order;
# every old value is unique in the document, so wdReplaceAll (2) only
# ever touches the single intended occurrence. Pairs are applied in the
# same order as the source diff so that a newly written value which
# happens to equal an earlier 'old' value never gets re-matched.

$d = $word.ActiveDocument

$d.Content.Find.Execute('51÷2=25, 1', $true, $false, $false, $false, $false, $true, 1, $false, '68÷2=34, 0', 2) | Out-Null
$d.Content.Find.Execute('87÷9=9, 6', $true, $false, $false, $false, $false, $true, 1, $false, '17÷4=4, 1', 2) | Out-Null
$d.Content.Find.Execute('16÷9=1, 7', $true, $false, $false, $false, $false, $true, 1, $false, '46÷9=5, 1', 2) | Out-Null
$d.Content.Find.Execute('36÷8=4, 4', $true, $false, $false, $false, $false, $true, 1, $false, '82÷7=11, 5', 2) | Out-Null
$d.Content.Find.Execute('25÷6=4, 1', $true, $false, $false, $false, $false, $true, 1, $false, '83÷4=20, 3', 2) | Out-Null
$d.Content.Find.Execute('10÷8=1, 2', $true, $false, $false, $false, $false, $true, 1, $false, '54÷6=9, 0', 2) | Out-Null
$d.Content.Find.Execute('18÷5=3, 3', $true, $false, $false, $false, $false, $true, 1, $false, '67÷2=33, 1', 2) | Out-Null
$d.Content.Find.Execute('96÷5=19, 1', $true, $false, $false, $false, $false, $true, 1, $false, '18÷2=9, 0', 2) | Out-Null
$d.Content.Find.Execute('38÷7=5, 3', $true, $false, $false, $false, $false, $true, 1, $false, '19÷5=3, 4', 2) | Out-Null
$d.Content.Find.Execute('13÷9=1, 4', $true, $false, $false, $false, $false, $true, 1, $false, '65÷5=13, 0', 2) | Out-Null
$d.Content.Find.Execute('23÷5=4, 3', $true, $false, $false, $false, $false, $true, 1, $false, '30÷3=10, 0', 2) | Out-Null
$d.Content.Find.Execute('82÷3=27, 1', $true, $false, $false, $false, $false, $true, 1, $false, '76÷6=12, 4', 2) | Out-Null
$d.Content.Find.Execute('17÷5=3, 2', $true, $false, $false, $false, $false, $true, 1, $false, '93÷3=31, 0', 2) | Out-Null
$d.Content.Find.Execute('31÷3=10, 1', $true, $false, $false, $false, $false, $true, 1, $false, '69÷6=11, 3', 2) | Out-Null
$d.Content.Find.Execute('52÷5=10, 2', $true, $false, $false, $false, $false, $true, 1, $false, '17÷4=4, 1', 2) | Out-Null
$d.Content.Find.Execute('58÷5=11, 3', $true, $false, $false, $false, $false, $true, 1, $false, '16÷9=1, 7', 2) | Out-Null
$d.Content.Find.Execute('15÷6=2, 3', $true, $false, $false, $false, $false, $true, 1, $false, '45÷8=5, 5', 2) | Out-Null
$d.Content.Find.Execute('90÷5=18, 0', $true, $false, $false, $false, $false, $true, 1, $false, '50÷3=16, 2', 2) | Out-Null
$d.Content.Find.Execute('64÷6=10, 4', $true, $false, $false, $false, $false, $true, 1, $false, '59÷4=14, 3', 2) | Out-Null
$d.Content.Find.Execute('31÷7=4, 3', $true, $false, $false, $false, $false, $true, 1, $false, '88÷7=12, 4', 2) | Out-Null
$d.Content.Find.Execute('12÷6=2, 0', $true, $false, $false, $false, $false, $true, 1, $false, '33÷5=6, 3', 2) | Out-Null
$d.Content.Find.Execute('49÷7=7, 0', $true, $false, $false, $false, $false, $true, 1, $false, '39÷3=13, 0', 2) | Out-Null
$d.Content.Find.Execute('44÷4=11, 0', $true, $false, $false, $false, $false, $true, 1, $false, '15÷2=7, 1', 2) | Out-Null
$d.Content.Find.Execute('73÷8=9, 1', $true, $false, $false, $false, $false, $true, 1, $false, '82÷4=20, 2', 2) | Out-Null
$d.Content.Find.Execute('49÷4=12, 1', $true, $false, $false, $false, $false, $true, 1, $false, '38÷2=19, 0', 2) | Out-Null

Write-Output "Replacements applied."
